$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that must stay text (avoid Excel auto-numeric conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.227.81"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.83"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.70"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4847"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3826"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07383"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9400"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.90"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07817"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.929.99"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.512"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.652"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.33"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008849"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.243.90"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.88"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.166"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.162.73"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.60"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.922"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.59"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.111"
$ws.Range("E28").Value = "  +4.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.42"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.976"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08910"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.255"
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7752"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.644"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02055"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05329"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5535"
$ws.Range("E40").Value = "  +2.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.005"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.071"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1533"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.476"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4866"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.13"
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.661"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.80"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06114"
$ws.Range("E51").Value = "  +0.07%  "
